$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Main sheet (covid19_cases_switzerland): fill in new daily figures
# ------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("covid19_cases_switzerland")

# Row 13 gets a new TG (col U) figure
$wsMain.Cells.Item(13, 21).Value = 23

# Row 14 is a brand-new day (2020-03-18) with partial cantonal updates
$wsMain.Cells.Item(14, 1).Value = 43908
$wsMain.Cells.Item(14, 1).NumberFormat = "yyyy\-mm\-dd;@"
$wsMain.Cells.Item(14, 5).Value = 193   # E  BE
$wsMain.Cells.Item(14, 7).Value = 182   # G  BS
$wsMain.Cells.Item(14, 14).Value = 99   # N  NE
$wsMain.Cells.Item(14, 15).Value = 12   # O  NW
$wsMain.Cells.Item(14, 17).Value = 61   # Q  SG
$wsMain.Cells.Item(14, 21).Value = 32   # U  TG
$wsMain.Cells.Item(14, 27).Value = 424  # AA ZH

# ------------------------------------------------------------------
# 2. Add new "Tests" worksheet right after the main sheet
# ------------------------------------------------------------------
$wsTests = $wb.Worksheets.Add($null, $wsMain)
$wsTests.Name = "Tests"

$wsTests.Columns.Item(1).ColumnWidth = 9.43
$wsTests.Columns.Item(2).ColumnWidth = 8.14

$headers = @("Date","AG","AI","AR","BE","BL","BS","FR","GE","GL","GR","JU","LU","NE","NW","OW","SG","SH","SO","SZ","TG","TI","UR","VD","VS","ZG","ZH","CH")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $wsTests.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$dates = @(43896,43897,43898,43899,43900,43901,43902,43903,43904,43905,43906,43907)
for ($r = 0; $r -lt $dates.Count; $r++) {
    $row = $r + 2
    $wsTests.Cells.Item($row, 1).Value = $dates[$r]
    $wsTests.Cells.Item($row, 1).NumberFormat = "yyyy\-mm\-dd;@"
    $wsTests.Cells.Item($row, 28).NumberFormat = "0"
}

# TG test count on the last available day
$wsTests.Cells.Item(13, 21).Value = 276

$wsTests.Range("U14").Select()

# ------------------------------------------------------------------
# 3. "Quellen" sheet: add TG source row + hyperlink
# ------------------------------------------------------------------
$wsQuellen = $wb.Worksheets.Item("Quellen")
$wsQuellen.Range("A17").Value = "TG"
$wsQuellen.Range("B17").Value = "https://www.tg.ch/news/fachdossier-coronavirus.html/10552"
$wsQuellen.Hyperlinks.Add($wsQuellen.Range("B17"), "https://www.tg.ch/news/fachdossier-coronavirus.html/10552")
$wsQuellen.Range("B17").Style = "Hyperlink"

$wsQuellen.Range("B2").Select()

# ------------------------------------------------------------------
# 4. Restore selection / active sheet state
# ------------------------------------------------------------------
$wsMain.Activate()
$wsMain.Range("AA15").Select()
